$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell 2 4 '23.694.14'
Set-TextCell 2 5 '  +1.14%  '
Set-TextCell 3 4 '1.650.42'
Set-TextCell 3 5 '  +1.20%  '
Set-TextCell 4 4 '0.9973'
Set-TextCell 4 5 '  -0.35%  '
Set-TextCell 5 4 '0.9979'
Set-TextCell 6 4 '304.47'
Set-TextCell 6 5 '  +0.04%  '
Set-TextCell 7 4 '0.3808'
Set-TextCell 7 5 '  +0.68%  '
Set-TextCell 8 4 '52.02'
Set-TextCell 8 5 '  +0.92%  '
Set-TextCell 9 4 '0.3613'
Set-TextCell 9 5 '  -1.02%  '
Set-TextCell 10 4 '1.249'
Set-TextCell 10 5 '  +1.24%  '
Set-TextCell 11 4 '0.08203'
Set-TextCell 11 5 '  -0.24%  '
Set-TextCell 12 4 '0.9987'
Set-TextCell 12 5 '  -0.24%  '
Set-TextCell 13 4 '22.57'
Set-TextCell 13 5 '  +0.94%  '
Set-TextCell 14 4 '6.538'
Set-TextCell 14 5 '  -0.27%  '
Set-TextCell 15 4 '7.393'
Set-TextCell 15 5 '  +0.68%  '
Set-TextCell 16 4 '0.00001232'
Set-TextCell 16 5 '  -1.64%  '
Set-TextCell 17 4 '1.649.41'
Set-TextCell 17 5 '  +1.35%  '
Set-TextCell 18 4 '96.91'
Set-TextCell 18 5 '  +3.03%  '
Set-TextCell 19 4 '0.06971'
Set-TextCell 19 5 '  -0.26%  '
Set-TextCell 20 4 '6.731'
Set-TextCell 20 5 '  +3.14%  '
Set-TextCell 21 4 '17.65'
Set-TextCell 21 5 '  -0.28%  '
Set-TextCell 22 4 '0.9984'
Set-TextCell 22 5 '  -0.26%  '
Set-TextCell 23 4 '12.58'
Set-TextCell 23 5 '  -1.18%  '
Set-TextCell 24 4 '23.646.33'
Set-TextCell 24 5 '  +0.99%  '
Set-TextCell 25 4 '2.519'
Set-TextCell 25 5 '  +2.60%  '
Set-TextCell 26 4 '3.108'
Set-TextCell 26 5 '  -1.61%  '
Set-TextCell 27 4 '21.30'
Set-TextCell 27 5 '  -0.47%  '
Set-TextCell 28 5 '  +1.14%  '
Set-TextCell 29 4 '5.195'
Set-TextCell 29 5 '  -2.06%  '
Set-TextCell 30 4 '134.81'
Set-TextCell 30 5 '  +0.56%  '
Set-TextCell 31 4 '1.834.98'
Set-TextCell 31 5 '  +1.45%  '
Set-TextCell 32 4 '6.812'
Set-TextCell 32 5 '  -0.41%  '
Set-TextCell 33 4 '1.093'
Set-TextCell 33 5 '  +6.92%  '
Set-TextCell 34 2 'WEMIXTOKEN'
Set-TextCell 34 3 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 34 4 '2.057'
Set-TextCell 34 5 '  -9.00%  '
Set-TextCell 35 2 'FraxShare'
Set-TextCell 35 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 35 4 '11.53'
Set-TextCell 35 5 '  +5.69%  '
Set-TextCell 36 4 '0.02814'
Set-TextCell 36 5 '  +0.93%  '
Set-TextCell 37 4 '0.2521'
Set-TextCell 37 5 '  -0.10%  '
Set-TextCell 38 4 '0.08829'
Set-TextCell 38 5 '  +0.65%  '
Set-TextCell 39 4 '6.089'
Set-TextCell 39 5 '  +0.71%  '
Set-TextCell 40 4 '0.07051'
Set-TextCell 40 5 '  -0.99%  '
Set-TextCell 41 4 '12.88'
Set-TextCell 41 5 '  +5.91%  '
Set-TextCell 42 4 '0.7070'
Set-TextCell 42 5 '  +0.42%  '
Set-TextCell 43 4 '1.332'
Set-TextCell 43 5 '  -1.34%  '
Set-TextCell 44 4 '15.88'
Set-TextCell 44 5 '  -2.08%  '
Set-TextCell 45 4 '0.6517'
Set-TextCell 45 5 '  -0.71%  '
Set-TextCell 46 4 '2.340'
Set-TextCell 46 5 '  +1.12%  '
Set-TextCell 47 4 '0.9984'
Set-TextCell 47 5 '  -0.18%  '
Set-TextCell 48 4 '3.980'
Set-TextCell 48 5 '  +0.11%  '
Set-TextCell 49 4 '0.07985'
Set-TextCell 49 5 '  -0.33%  '
Set-TextCell 50 4 '128.12'
Set-TextCell 50 5 '  +1.43%  '
Set-TextCell 51 4 '1.196'
Set-TextCell 51 5 '  -0.18%  '
